# Add a TOTAL row (row 4) to the "CUMPLIMIENTO MENSUAL" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# B4: "TOTAL" label, right-aligned (new style: numFmtId 0 + right alignment)
$ws.Range("B4").Value = "TOTAL"
$ws.Range("B4").HorizontalAlignment = -4152  # xlRight

# C4:E4: currency totals - reuse the existing currency format from row 3
$ws.Range("C3:E3").Copy()
$ws.Range("C4:E4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C4").Value = 17500
$ws.Range("D4").Value = 605.48
$ws.Range("E4").Value = 16894.52

# F4: percentage cumplimiento - reuse the existing percent format from row 3
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F4").Value = 0.03459885714285715
